$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Update of the diary" -----------------------------------------------------

# Tweak wording of the section header above the attendance grid.
$ws.Range("C3").Value = "Presenza membri del gruppo"

# The "Sistemi e reti" diary entry becomes just "Sistemi".
$ws.Range("B10").Value = "Sistemi"

# Append the new diary entries (classes held on 16/05 and 17/05).
$ws.Range("A11").Value = 45428
$ws.Range("B11").Value = "Storia"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0

$ws.Range("A12").Value = 45429
$ws.Range("B12").Value = "Informatica 1"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1

$ws.Range("A13").Value = 45429
$ws.Range("B13").Value = "Informatica 2"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1

# Carry the same cell formatting used by the existing diary rows onto the
# three new rows (xlPasteFormats = -4122).
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F13").PasteSpecial(-4122)

# The later diary entries use the tighter 13.8pt row height.
$ws.Range("A5:F5").RowHeight = 13.8
$ws.Range("A11:F11").RowHeight = 13.8
$ws.Range("A12:F12").RowHeight = 13.8
$ws.Range("A13:F13").RowHeight = 13.8

# Leave the selection where the author left it before saving.
$ws.Range("G2").Select()

$wb.Application.CalculateFull()
